# Applies the "Automatic update of files." diff to the "Avverkningsanmälningar"
# worksheet:
#   1. The "Förändrad" column (C) is bumped by one day (46077 -> 46078) for every
#      data row (rows 2..54).
#   2. Rows 5..54 are re-ordered (the underlying records are the same, only their
#      row position changes) according to the mapping baked into $oldRows/$newRows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 26   # column Z

# New row position (index) -> original row that should end up there.
$newRows = @(5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51,52,53,54)
$oldRows = @(6,7,8,5,9,12,10,11,13,14,15,37,38,17,41,46,20,35,44,21,16,49,43,54,22,23,42,45,27,30,39,28,47,50,51,52,48,53,29,36,25,31,40,26,18,32,24,33,34,19)

# 1. Snapshot every source row (5..54) completely before any writes happen, so
#    that overlapping source/destination rows do not clobber data we still need.
$snapshots = @{}
foreach ($r in $oldRows) {
    if (-not $snapshots.ContainsKey($r)) {
        $rowData = @{}
        for ($c = 1; $c -le $lastCol; $c++) {
            $rowData[$c] = $ws.Cells.Item($r, $c).Formula
        }
        $snapshots[$r] = $rowData
    }
}

# 2. Write each snapshot back into its new row position.
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $nr = $newRows[$i]
    $or = $oldRows[$i]
    $rowData = $snapshots[$or]
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($nr, $c).Formula = $rowData[$c]
    }
}

# 3. Bump the "Förändrad" (column C) date by one day for every data row (2..54).
for ($r = 2; $r -le 54; $r++) {
    $cur = $ws.Cells.Item($r, 3).Value2
    if ($cur -ne $null) {
        $ws.Cells.Item($r, 3).Value2 = $cur + 1
    }
}

# 4. Re-writing a wrapped multi-line cell (column R) makes the host engine
#    auto-fit the row height. The source workbook instead keeps an explicit
#    15pt row height for every populated row (row 54 is the lone exception,
#    left at the sheet default), so restore that here.
for ($r = 5; $r -le 53; $r++) {
    $ws.Rows.Item($r).RowHeight = 15
}
